# Daily attendance processing - 2025-11-07 18:29:07
# Normalizes the "Recorded By" (column G) cell values so that the
# dnasr281@gmail.com (or admin@admin.com) reviewer account is listed
# before the "System" account that originally recorded the session,
# and the lowercase "system" backup alias is listed before the
# capitalized "System" entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($i = 2; $i -le $lastRow; $i++) {
    $cell = $ws.Cells.Item($i, 7)
    $val = $cell.Value2

    if ($val -eq $null) {
        continue
    }

    $parts = $val -split ", "

    if ($parts.Count -eq 2 -and $parts[1] -eq "dnasr281@gmail.com") {
        $cell.Value2 = $parts[1] + ", " + $parts[0]
    }
    elseif ($parts.Count -eq 2 -and $parts[0] -eq "admin@admin.com" -and $parts[1] -eq "dnasr281@gmail.com") {
        $cell.Value2 = $parts[1] + ", " + $parts[0]
    }
    elseif ($parts.Count -eq 3 -and $parts[0] -eq "System" -and $parts[1] -eq "system") {
        $cell.Value2 = $parts[1] + ", " + $parts[0] + ", " + $parts[2]
    }
}
